$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mistranslated / placeholder header values on row 2 (G2:I2) with the
# correct French labels used by the canteens import fixture.
$ws.Range("G2").Value = "Restaurant avec cuisine sur place"
$ws.Range("H2").Value = "Concédée"
$ws.Range("I2").Value = "Public"

# G2 ("type_production") switches to a Times New Roman font while keeping a
# text number format.
$ws.Range("G2").Font.Name = "Times New Roman"
$ws.Range("G2").NumberFormat = "@"

# H2 / I2 ("type_gestion", "modèle_économique") switch to a wrapped,
# General-formatted cell instead of plain text.
$ws.Range("H2:I2").NumberFormat = "General"
$ws.Range("H2:I2").WrapText = $true

# Move the active selection from F2 to I2.
$ws.Range("I2").Select()
